$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value2 = 118.6930263333333
$ws.Range("H2").Value2 = 356.079079
$ws.Range("I2").Value2 = 0.2696481350657977
$ws.Range("J2").Value2 = 0.2696481350657977
$ws.Range("K2").Value2 = 2
$ws.Range("L2").Value2 = 0.6666666666666666
$ws.Range("M2").Value2 = 1.417175333333333
$ws.Range("N2").Value2 = 4.251526
$ws.Range("O2").Value2 = 0.02388929759887871
$ws.Range("P2").Value2 = 0.02388929759887871
$ws.Range("Q2").Value2 = 168.2088291582838
$ws.Range("R2").Value2 = 1513.879462424554
$ws.Range("S2").Value2 = 0.006441704545569481
$ws.Range("T2").Value2 = 0.006441704545569481
# Row 3
$ws.Range("G3").Value2 = 118.6930263333333
$ws.Range("H3").Value2 = 356.079079
$ws.Range("I3").Value2 = 0.2696481350657977
$ws.Range("J3").Value2 = 0.2696481350657977
$ws.Range("O3").Value2 = 0.0375290450698413
$ws.Range("P3").Value2 = 0.03752904506984131
$ws.Range("Q3").Value2 = 264.2487375151109
$ws.Range("R3").Value2 = 2378.238637635998
$ws.Range("S3").Value2 = 0.01011963701388298
$ws.Range("T3").Value2 = 0.01011963701388298
# Row 4
$ws.Range("G4").Value2 = 118.6930263333333
$ws.Range("H4").Value2 = 356.079079
$ws.Range("I4").Value2 = 0.2696481350657977
$ws.Range("J4").Value2 = 0.2696481350657977
$ws.Range("M4").Value2 = 55.67910766666667
$ws.Range("N4").Value2 = 167.037323
$ws.Range("O4").Value2 = 0.93858165733128
$ws.Range("P4").Value2 = 0.93858165733128
$ws.Range("Q4").Value2 = 6608.721792496169
$ws.Range("R4").Value2 = 59478.49613246552
$ws.Range("S4").Value2 = 0.2530867935063452
$ws.Range("T4").Value2 = 0.2530867935063452
# Row 5
$ws.Range("I5").Value2 = 0.3482063679522526
$ws.Range("J5").Value2 = 0.3482063679522526
$ws.Range("K5").Value2 = 2
$ws.Range("L5").Value2 = 0.6666666666666666
$ws.Range("M5").Value2 = 1.417175333333333
$ws.Range("N5").Value2 = 4.251526
$ws.Range("O5").Value2 = 0.02388929759887871
$ws.Range("P5").Value2 = 0.02388929759887871
$ws.Range("Q5").Value2 = 217.2141314621545
$ws.Range("R5").Value2 = 1954.92718315939
$ws.Range("S5").Value2 = 0.008318405549836025
$ws.Range("T5").Value2 = 0.008318405549836025
# Row 6
$ws.Range("I6").Value2 = 0.3482063679522526
$ws.Range("J6").Value2 = 0.3482063679522526
$ws.Range("O6").Value2 = 0.0375290450698413
$ws.Range("P6").Value2 = 0.03752904506984131
$ws.Range("S6").Value2 = 0.01306785247648583
$ws.Range("T6").Value2 = 0.01306785247648584
# Row 7
$ws.Range("I7").Value2 = 0.3482063679522526
$ws.Range("J7").Value2 = 0.3482063679522526
$ws.Range("M7").Value2 = 55.67910766666667
$ws.Range("N7").Value2 = 167.037323
$ws.Range("O7").Value2 = 0.93858165733128
$ws.Range("P7").Value2 = 0.93858165733128
$ws.Range("Q7").Value2 = 8534.080948160346
$ws.Range("R7").Value2 = 76806.72853344311
$ws.Range("S7").Value2 = 0.3268201099259307
$ws.Range("T7").Value2 = 0.3268201099259307
# Row 8
$ws.Range("G8").Value2 = 116.0670876666667
$ws.Range("H8").Value2 = 348.201263
$ws.Range("I8").Value2 = 0.2636824984472209
$ws.Range("J8").Value2 = 0.2636824984472209
$ws.Range("K8").Value2 = 2
$ws.Range("L8").Value2 = 0.6666666666666666
$ws.Range("M8").Value2 = 1.417175333333333
$ws.Range("N8").Value2 = 4.251526
$ws.Range("O8").Value2 = 0.02388929759887871
$ws.Range("P8").Value2 = 0.02388929759887871
$ws.Range("Q8").Value2 = 164.4874136530376
$ws.Range("R8").Value2 = 1480.386722877338
$ws.Range("S8").Value2 = 0.006299189677021533
$ws.Range("T8").Value2 = 0.006299189677021533
# Row 9
$ws.Range("G9").Value2 = 116.0670876666667
$ws.Range("H9").Value2 = 348.201263
$ws.Range("I9").Value2 = 0.2636824984472209
$ws.Range("J9").Value2 = 0.2636824984472209
$ws.Range("O9").Value2 = 0.0375290450698413
$ws.Range("P9").Value2 = 0.03752904506984131
$ws.Range("Q9").Value2 = 258.4025559921118
$ws.Range("R9").Value2 = 2325.623003929006
$ws.Range("S9").Value2 = 0.009895752368354114
$ws.Range("T9").Value2 = 0.009895752368354115
# Row 10
$ws.Range("G10").Value2 = 116.0670876666667
$ws.Range("H10").Value2 = 348.201263
$ws.Range("I10").Value2 = 0.2636824984472209
$ws.Range("J10").Value2 = 0.2636824984472209
$ws.Range("M10").Value2 = 55.67910766666667
$ws.Range("N10").Value2 = 167.037323
$ws.Range("O10").Value2 = 0.93858165733128
$ws.Range("P10").Value2 = 0.93858165733128
$ws.Range("Q10").Value2 = 6462.511870748774
$ws.Range("R10").Value2 = 58162.60683673896
$ws.Range("S10").Value2 = 0.2474875564018452
$ws.Range("T10").Value2 = 0.2474875564018452
# Row 11
$ws.Range("G11").Value2 = 52.14473966666667
$ws.Range("H11").Value2 = 156.434219
$ws.Range("I11").Value2 = 0.1184629985347288
$ws.Range("J11").Value2 = 0.1184629985347288
$ws.Range("K11").Value2 = 2
$ws.Range("L11").Value2 = 0.6666666666666666
$ws.Range("M11").Value2 = 1.417175333333333
$ws.Range("N11").Value2 = 4.251526
$ws.Range("O11").Value2 = 0.02388929759887871
$ws.Range("P11").Value2 = 0.02388929759887871
$ws.Range("Q11").Value2 = 73.89823881868824
$ws.Range("R11").Value2 = 665.0841493681941
$ws.Range("S11").Value2 = 0.002829997826451668
$ws.Range("T11").Value2 = 0.002829997826451669
# Row 12
$ws.Range("G12").Value2 = 52.14473966666667
$ws.Range("H12").Value2 = 156.434219
$ws.Range("I12").Value2 = 0.1184629985347288
$ws.Range("J12").Value2 = 0.1184629985347288
$ws.Range("O12").Value2 = 0.0375290450698413
$ws.Range("P12").Value2 = 0.03752904506984131
$ws.Range("Q12").Value2 = 116.0909115778531
$ws.Range("R12").Value2 = 1044.818204200678
$ws.Range("S12").Value2 = 0.004445803211118381
$ws.Range("T12").Value2 = 0.004445803211118382
# Row 13
$ws.Range("G13").Value2 = 52.14473966666667
$ws.Range("H13").Value2 = 156.434219
$ws.Range("I13").Value2 = 0.1184629985347288
$ws.Range("J13").Value2 = 0.1184629985347288
$ws.Range("M13").Value2 = 55.67910766666667
$ws.Range("N13").Value2 = 167.037323
$ws.Range("O13").Value2 = 0.93858165733128
$ws.Range("P13").Value2 = 0.93858165733128
$ws.Range("Q13").Value2 = 2903.372574150638
$ws.Range("R13").Value2 = 26130.35316735574
$ws.Range("S13").Value2 = 0.1111871974971587
$ws.Range("T13").Value2 = 0.1111871974971587
